# Applies odds-value corrections to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.75
$ws.Range("I5").Value = 2.75
$ws.Range("AL5").Value = 26
# Row 10
$ws.Range("N10").Value = 10
# Row 13
$ws.Range("H13").Value = 3.65
$ws.Range("J13").Value = 2.92
$ws.Range("K13").Value = 2.35
$ws.Range("L13").Value = 2.9
$ws.Range("U13").Value = 1.37
$ws.Range("V13").Value = 2.85
$ws.Range("W13").Value = 15.5
$ws.Range("X13").Value = 18.5
$ws.Range("Z13").Value = 32
$ws.Range("AB13").Value = 18
$ws.Range("AD13").Value = 8
$ws.Range("AE13").Value = 10.25
$ws.Range("AH13").Value = 14
$ws.Range("AI13").Value = 17
$ws.Range("AK13").Value = 29
$ws.Range("AL13").Value = 17
$ws.Range("AM13").Value = 18
$ws.Range("AN13").Value = 5.1
$ws.Range("AO13").Value = 12.5
$ws.Range("AP13").Value = 15
$ws.Range("AQ13").Value = 45
$ws.Range("AR13").Value = 55
$ws.Range("AS13").Value = 120
$ws.Range("AX13").Value = 12.5
$ws.Range("AY13").Value = 15.5
$ws.Range("AZ13").Value = 45
$ws.Range("BA13").Value = 55
# Row 18
$ws.Range("K18").Value = 2.4
$ws.Range("L18").Value = 2.05
$ws.Range("M18").Value = 1.03
$ws.Range("O18").Value = 1.18
$ws.Range("S18").Value = 1.3
$ws.Range("T18").Value = 3.4
$ws.Range("AT18").Value = 3.4
# Row 19
$ws.Range("M19").Value = 1.04
$ws.Range("O19").Value = 1.2
# Row 20
$ws.Range("J20").Value = 2.6
$ws.Range("L20").Value = 4
$ws.Range("Q20").Value = 1.75
$ws.Range("R20").Value = 2.05
$ws.Range("AN20").Value = 4
$ws.Range("AO20").Value = 10
$ws.Range("AP20").Value = 19
$ws.Range("AQ20").Value = 34
$ws.Range("AT20").Value = 3.25
$ws.Range("AW20").Value = 5.5
$ws.Range("AX20").Value = 19
# Row 24
$ws.Range("J24").Value = 3.25
$ws.Range("M24").Value = 1.11
$ws.Range("N24").Value = 6.5
$ws.Range("AN24").Value = 4.33
$ws.Range("AR24").Value = 101
$ws.Range("AY24").Value = 34
$ws.Range("AZ24").Value = 67
# Row 28
$ws.Range("N28").Value = 15
# Row 29
$ws.Range("G29").Value = 1.22
$ws.Range("H29").Value = 6.25
$ws.Range("I29").Value = 13
$ws.Range("J29").Value = 1.62
$ws.Range("K29").Value = 2.88
$ws.Range("L29").Value = 9.5
$ws.Range("Q29").Value = 1.48
$ws.Range("R29").Value = 2.6
$ws.Range("AB29").Value = 29
$ws.Range("AC29").Value = 17
$ws.Range("AD29").Value = 12
$ws.Range("AJ29").Value = 34
$ws.Range("AL29").Value = 81
$ws.Range("AX29").Value = 51
$ws.Range("AZ29").Value = 251
# Row 30
$ws.Range("G30").Value = 2.55
$ws.Range("I30").Value = 2.55
$ws.Range("L30").Value = 2.88
$ws.Range("Q30").Value = 1.4
$ws.Range("R30").Value = 2.88
$ws.Range("U30").Value = 1.36
$ws.Range("V30").Value = 3
$ws.Range("Y30").Value = 11
$ws.Range("AD30").Value = 8
$ws.Range("AH30").Value = 15
$ws.Range("AI30").Value = 17
$ws.Range("AK30").Value = 26
$ws.Range("AL30").Value = 17
$ws.Range("AM30").Value = 19
$ws.Range("BC30").Value = 201
# Row 31
$ws.Range("O31").Value = 1.13
$ws.Range("P31").Value = 6
# Row 35
$ws.Range("G35").Value = 7
$ws.Range("H35").Value = 5.25
$ws.Range("I35").Value = 1.33
$ws.Range("L35").Value = 1.8
$ws.Range("AD35").Value = 11
$ws.Range("AE35").Value = 19
$ws.Range("AG35").Value = 201
$ws.Range("AJ35").Value = 9
$ws.Range("AK35").Value = 9.5
$ws.Range("AN35").Value = 9
$ws.Range("AQ35").Value = 126
$ws.Range("AR35").Value = 126
$ws.Range("AS35").Value = 201
$ws.Range("AU35").Value = 8.5
$ws.Range("AV35").Value = 51
# Row 39
$ws.Range("J39").Value = 7
$ws.Range("O39").Value = 1.22
$ws.Range("P39").Value = 4
$ws.Range("S39").Value = 1.33
$ws.Range("T39").Value = 3.25
$ws.Range("U39").Value = 2
$ws.Range("V39").Value = 1.73
$ws.Range("Z39").Value = 81
$ws.Range("AG39").Value = 351
$ws.Range("AH39").Value = 7
$ws.Range("AS39").Value = 301
$ws.Range("AT39").Value = 3.25
$ws.Range("AU39").Value = 9
$ws.Range("AV39").Value = 51
$ws.Range("AW39").Value = 3.4
$ws.Range("AX39").Value = 7
# Row 41
$ws.Range("AY41").Value = 29
$ws.Range("BB41").Value = 201
# Row 42
$ws.Range("M42").Value = 1.08
$ws.Range("N42").Value = 8
$ws.Range("Q42").Value = 2.3
$ws.Range("R42").Value = 1.6
# Row 49
$ws.Range("O49").Value = 1.17
$ws.Range("P49").Value = 5
# Row 57
$ws.Range("G57").Value = 1.48
$ws.Range("I57").Value = 7
$ws.Range("O57").Value = 1.2
$ws.Range("P57").Value = 4.33
$ws.Range("Q57").Value = 1.7
$ws.Range("R57").Value = 2.1
$ws.Range("AM57").Value = 41
$ws.Range("AZ57").Value = 101
# Row 77
$ws.Range("Q77").Value = 1.8
# Row 78
$ws.Range("R78").Value = 1.7
# Row 82
$ws.Range("I82").Value = 4.2
$ws.Range("J82").Value = 2.88
$ws.Range("K82").Value = 1.91
$ws.Range("L82").Value = 5
$ws.Range("M82").Value = 1.13
$ws.Range("N82").Value = 6
$ws.Range("O82").Value = 1.53
$ws.Range("P82").Value = 2.38
$ws.Range("Q82").Value = 2.7
$ws.Range("R82").Value = 1.44
$ws.Range("S82").Value = 1.62
$ws.Range("T82").Value = 2.2
$ws.Range("U82").Value = 2.38
$ws.Range("V82").Value = 1.53
$ws.Range("W82").Value = 5
$ws.Range("Y82").Value = 10
$ws.Range("AC82").Value = 6
$ws.Range("AD82").Value = 6.5
$ws.Range("AH82").Value = 8.5
$ws.Range("AK82").Value = 51
$ws.Range("AS82").Value = 301
$ws.Range("AT82").Value = 2.2
$ws.Range("AU82").Value = 10
$ws.Range("AW82").Value = 6
$ws.Range("AZ82").Value = 101
# Row 90
$ws.Range("N90").Value = 13
$ws.Range("Q90").Value = 1.85
$ws.Range("R90").Value = 2
# Row 91
$ws.Range("M91").Value = 1.07
$ws.Range("N91").Value = 9
$ws.Range("Q91").Value = 2.1
$ws.Range("R91").Value = 1.7
# Row 96
$ws.Range("G96").Value = 9.5
$ws.Range("N96").Value = 12
$ws.Range("O96").Value = 1.26
$ws.Range("P96").Value = 3.2
$ws.Range("Q96").Value = 1.82
$ws.Range("R96").Value = 1.9
$ws.Range("S96").Value = 1.36
$ws.Range("T96").Value = 3
$ws.Range("U96").Value = 2.27
$ws.Range("W96").Value = 20
$ws.Range("X96").Value = 65
$ws.Range("AD96").Value = 9.75
$ws.Range("AH96").Value = 5.8
$ws.Range("AI96").Value = 5.4
$ws.Range("AJ96").Value = 9.25
$ws.Range("AL96").Value = 12.5
$ws.Range("AM96").Value = 40
$ws.Range("AO96").Value = 60
$ws.Range("AX96").Value = 5.6
$ws.Range("AY96").Value = 19
$ws.Range("BB96").Value = 350
